# Fix text wrapping and cell merging in extended Executive Summary sections
# for full text visibility.
#
# This script:
#   1. Resizes the "Executive Summary" columns so column A is wide and
#      B:E become narrow spacer columns.
#   2. Restyles / re-heights / merges rows 17-20 (Trading Activity Summary),
#      23-28 (Key Insights) and 31-36 (Action Items) so each row's text
#      lives in a single wrapped, merged A:E cell instead of being split
#      across separate B/C/D/E value cells.
#   3. Clears the now-redundant B:E text that used to hold those values.
#   4. Mirrors the same "blank divider" styling (already used elsewhere in
#      the workbook) onto the Monthly Performance sheet's 12-MONTH SUMMARY
#      header row (row 33).

$wb = $excel.ActiveWorkbook
$wsES = $wb.Worksheets.Item("Executive Summary")
$wsMP = $wb.Worksheets.Item("Monthly Performance")

# ---------------------------------------------------------------------
# 1. Column widths on the Executive Summary sheet
# ---------------------------------------------------------------------
# NOTE: the runtime's ColumnWidth setter stores width + 5/6 internally,
# so subtract 5/6 here in order to land on the exact target width.
$offset = 5 / 6
$wsES.Columns.Item(1).ColumnWidth = 95 - $offset
$wsES.Columns.Item(2).ColumnWidth = 2 - $offset
$wsES.Columns.Item(3).ColumnWidth = 2 - $offset
$wsES.Columns.Item(4).ColumnWidth = 2 - $offset
$wsES.Columns.Item(5).ColumnWidth = 2 - $offset

# ---------------------------------------------------------------------
# 2. TRADING ACTIVITY SUMMARY block (rows 17-20)
# ---------------------------------------------------------------------
# Grab an existing "left/center/wrap" style that already shares the same
# font as the target look (font used by the Key Insights labels), then
# recolor the fill to match the Trading Activity Summary section so the
# new combined style reuses as many existing resources as possible.
$tradingRows = 17, 18, 19, 20
foreach ($rowNum in $tradingRows) {
    $aCell = $wsES.Cells.Item($rowNum, 1)

    $wsES.Range("A23").Copy()
    $aCell.PasteSpecial(-4122)
    $aCell.Interior.Color = 13431551   # FFF2CC (matches Trading Activity fill)
    $aCell.WrapText = $true
    $aCell.HorizontalAlignment = -4131  # xlLeft
    $aCell.VerticalAlignment = -4108    # xlCenter

    $wsES.Rows.Item($rowNum).RowHeight = 25

    # Clear the old per-column values - the row's message now lives
    # entirely in the merged A:E cell.
    $wsES.Range("B" + $rowNum + ":E" + $rowNum).ClearContents()

    # Apply the existing "blank divider" styles (same ones used
    # elsewhere in the workbook) to the B:D / E filler cells.
    $wsMP.Range("B28").Copy()
    $wsES.Range("B" + $rowNum + ":D" + $rowNum).PasteSpecial(-4122)
    $wsMP.Range("M28").Copy()
    $wsES.Range("E" + $rowNum).PasteSpecial(-4122)

    $wsES.Range("A" + $rowNum + ":E" + $rowNum).Merge()
}

# ---------------------------------------------------------------------
# 3. KEY INSIGHTS & RECOMMENDATIONS block (rows 23-28)
# ---------------------------------------------------------------------
$insightRows = 23, 24, 25, 26, 27, 28
foreach ($rowNum in $insightRows) {
    $aCell = $wsES.Cells.Item($rowNum, 1)

    # Keep the existing font/fill/border - only switch vertical alignment
    # from "center" to "top" so long wrapped text starts at the top of
    # the (now taller) row.
    $aCell.VerticalAlignment = -4160   # xlTop

    $wsES.Rows.Item($rowNum).RowHeight = 30

    # The second-column commentary text is no longer needed - the A-cell
    # now spans the full row.
    $wsES.Range("B" + $rowNum + ":E" + $rowNum).ClearContents()

    $wsMP.Range("B28").Copy()
    $wsES.Range("B" + $rowNum + ":D" + $rowNum).PasteSpecial(-4122)
    $wsMP.Range("M28").Copy()
    $wsES.Range("E" + $rowNum).PasteSpecial(-4122)

    $wsES.Range("A" + $rowNum + ":E" + $rowNum).Merge()
}

# ---------------------------------------------------------------------
# 4. ACTION ITEMS & STRATEGY block (rows 31-36)
# ---------------------------------------------------------------------
$actionRows = 31, 32, 33, 34, 35, 36
foreach ($rowNum in $actionRows) {
    $aCell = $wsES.Cells.Item($rowNum, 1)

    $aCell.VerticalAlignment = -4160   # xlTop

    $wsES.Rows.Item($rowNum).RowHeight = 30

    $wsES.Range("B" + $rowNum + ":E" + $rowNum).ClearContents()

    $wsMP.Range("B28").Copy()
    $wsES.Range("B" + $rowNum + ":D" + $rowNum).PasteSpecial(-4122)
    $wsMP.Range("M28").Copy()
    $wsES.Range("E" + $rowNum).PasteSpecial(-4122)

    $wsES.Range("A" + $rowNum + ":E" + $rowNum).Merge()
}

# ---------------------------------------------------------------------
# 5. Monthly Performance sheet - 12-MONTH SUMMARY header row (row 33)
# ---------------------------------------------------------------------
$wsMP.Range("B28").Copy()
$wsMP.Range("B33:L33").PasteSpecial(-4122)
$wsMP.Range("M28").Copy()
$wsMP.Range("M33").PasteSpecial(-4122)

$excel.CutCopyMode = $false
